# Rename the inline Pearson/BTEC logo pictures in the document's header
# and footer stories.
#
#   - Footer (default / primary)    -> id 2 : image1.png -> image2.png
#   - Footer (first page)           -> id 3 : image1.png -> image2.png
#   - Header (first page)           -> id 1 : image2.jpg -> image1.jpg
#
# `InlineShape` has no writable `.Name` of its own (matches real Word's
# object model), so each picture is promoted to a floating `Shape` long
# enough to rename it, then converted back to an inline shape in place.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($range, $newName) {
    $ishp = $range.InlineShapes.Item(1)
    $shp = $ishp.ConvertToShape()
    $shp.Name = $newName
    $shp.ConvertToInlineShape() | Out-Null
}

# Default (primary) footer picture: image1.png -> image2.png
Rename-InlinePicture $sec.Footers.Item(1).Range "image2.png"

# First-page footer picture: image1.png -> image2.png
Rename-InlinePicture $sec.Footers.Item(2).Range "image2.png"

# First-page header picture: image2.jpg -> image1.jpg
Rename-InlinePicture $sec.Headers.Item(2).Range "image1.jpg"

Write-Output "Renamed header/footer logo pictures."
